$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.498.93'
$ws.Range('E2').Value = '  +0.09%  '

$ws.Range('D3').Value = '1.565.40'
$ws.Range('E3').Value = '  -2.12%  '

$ws.Range('E4').Value = '  +0.16%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '211.71'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.54%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.492'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.23%  '

$ws.Range('E7').Value = '  +0.22%  '

$ws.Range('E8').Value = '  +4.69%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '24.02'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.43%  '

$ws.Range('E10').Value = '  -1.90%  '

$ws.Range('E11').Value = '  -1.64%  '

$ws.Range('E12').Value = '  -0.33%  '

$ws.Range('E13').Value = '  -2.05%  '

$ws.Range('D14').Value = '1.564.05'
$ws.Range('E14').Value = '  -2.17%  '

$ws.Range('E15').Value = '  -2.86%  '

$ws.Range('D16').Value = '28.521.02'
$ws.Range('E16').Value = '  +0.15%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.67'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -3.56%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '62.18'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -1.93%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '228.75'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.83%  '

$ws.Range('B20').Value = 'Chainlink'
$ws.Range('C20').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.34'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -2.63%  '

$ws.Range('B21').Value = 'ShibaInu'
$ws.Range('C21').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D21').Value = '0.0₃0693'
$ws.Range('E21').Value = '  -2.71%  '

$ws.Range('E22').Value = '  +0.00%  '

$ws.Range('E23').Value = '  -6.08%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '9.14'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -3.08%  '

$ws.Range('E25').Value = '  +6.52%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '150.65'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -1.46%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '14.97'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -2.25%  '

$ws.Range('E28').Value = '  -2.93%  '

$ws.Range('E29').Value = '  -4.21%  '

$ws.Range('E30').Value = '  +0.13%  '

$ws.Range('E31').Value = '  -2.16%  '

$ws.Range('E32').Value = '  -4.03%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.20'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -1.52%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.09'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -2.89%  '

$ws.Range('D35').Value = '1.390.20'
$ws.Range('E35').Value = '  -2.33%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.04'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.42%  '

$ws.Range('E37').Value = '  -3.47%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.35'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.95%  '

$ws.Range('E39').Value = '  +1.96%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0165'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.31%  '

$ws.Range('E41').Value = '  -1.81%  '

$ws.Range('E42').Value = '  +0.09%  '

$ws.Range('E43').Value = '  +2.83%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.788'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -4.26%  '

$ws.Range('E45').Value = '  -4.40%  '

$ws.Range('E46').Value = '  -0.38%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '62.70'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -3.54%  '

$ws.Range('D48').Value = '1.702.29'
$ws.Range('E48').Value = '  -1.99%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '86.12'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -1.73%  '

$ws.Range('E50').Value = '  -4.38%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0524'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.49%  '
